$d = $word.ActiveDocument
$d.Content.Find.Execute("203.0.113.24riot", $true, $false, $false, $false, $false,
                         $true, 1, $false, "203.0.113.24", 2)
